$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Enrollment_Race": update enrollment counts in column A for several
# rows (re-implemented enrollment-by-race numbers).
# ---------------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Enrollment_Race")
$wsRace.Range("A3").Value = 190
$wsRace.Range("A5").Value = 70
$wsRace.Range("A7").Value = 62
$wsRace.Range("A9").Value = 939
$wsRace.Range("A11").Value = 2
$wsRace.Range("A13").Value = 91
$wsRace.Range("A17").Value = 68
$wsRace.Range("A19").Value = 16

# ---------------------------------------------------------------------------
# Sheet "High School Units": re-implemented unit counts. The data for rows
# 6-21 effectively shift up by one (values that used to live in row N+1 now
# live in row N), row 22 gets new combined values, and the old trailing rows
# 23-24 are removed entirely.
# ---------------------------------------------------------------------------
$wsHS = $wb.Worksheets.Item("High School Units")

function Set-HSRow($row, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n) {
    $wsHS.Range("A$row").Value = $a
    $wsHS.Range("B$row").Value = $b
    $wsHS.Range("C$row").Value = $c
    $wsHS.Range("D$row").Value = $d
    $wsHS.Range("E$row").Value = $e
    $wsHS.Range("F$row").Value = $f
    $wsHS.Range("G$row").Value = $g
    $wsHS.Range("H$row").Value = $h
    $wsHS.Range("I$row").Value = $i
    $wsHS.Range("J$row").Value = $j
    $wsHS.Range("K$row").Value = $k
    $wsHS.Range("L$row").Value = $l
    $wsHS.Range("M$row").Value = $m
    $wsHS.Range("N$row").Value = $n
}

Set-HSRow 6  4 0 0 1 0 0 0 0 0 0 0 0 0 1
Set-HSRow 7  3 0 0 1 1 0 0 0 0 0 0 0 1 0
Set-HSRow 8  0 0 0 0 0 1 0 0 0 0 0 0 1 0
Set-HSRow 9  0 0 0 0 0 1 0 0 0 0 0 0 0 1
Set-HSRow 10 2 0 0 0 0 0 1 0 0 0 0 0 1 0
Set-HSRow 11 0 0 0 0 0 0 1 0 0 0 0 0 0 1
Set-HSRow 12 0 0 0 0 0 0 0 1 0 0 0 0 1 0
Set-HSRow 13 0 0 0 0 0 0 0 1 0 0 0 0 0 1
Set-HSRow 14 0 0 0 0 0 0 0 1 0 0 0 0 1 0
Set-HSRow 15 0 0 0 0 0 0 0 1 0 0 0 0 0 1
Set-HSRow 16 0 0 0 0 0 0 0 0 1 0 0 0 1 0
Set-HSRow 17 4 0 0 0 0 0 0 0 1 0 0 0 0 1
Set-HSRow 18 0 0 0 0 0 0 0 0 0 1 0 0 1 0
Set-HSRow 19 0 0 0 0 0 0 0 0 0 1 0 0 0 1
Set-HSRow 20 0 0 0 0 0 0 0 0 0 0 1 0 1 0
Set-HSRow 21 0 0 0 0 0 0 0 0 0 0 1 0 0 1
Set-HSRow 22 0 0 0 0 0 0 0 0 0 0 0 1 1 0

# Remove the now-obsolete trailing rows 23 and 24.
$wsHS.Rows("23:24").Delete()

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping. "High School Units" was the active
# tab and is now just a background sheet with a new selection; the newly
# active tab is "Enrollment_Race".
# ---------------------------------------------------------------------------
$wsHS.Activate()
$wsHS.Range("K15").Select()

$wsRace.Activate()
$wsRace.Range("M31").Select()
